$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "counter" values in column H (keep cell style/format)
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("H9").ClearContents()

# Update the active selection to G14
$ws.Range("G14").Select()
